$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 533 (weekly price update), shifting all
# following rows (533-557) down to (534-558).
$ws.Rows.Item(533).Insert()

$ws.Cells.Item(533, 1).Value = 4
$ws.Cells.Item(533, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(533, 3).Value = "Los Lagos"
$ws.Cells.Item(533, 4).Value = 45267
$ws.Cells.Item(533, 5).Value = 10
$ws.Cells.Item(533, 6).Value = 100112040
$ws.Cells.Item(533, 7).Value = "Cilantro"
$ws.Cells.Item(533, 8).Value = "Sin especificar"
$ws.Cells.Item(533, 9).Value = "Primera"
$ws.Cells.Item(533, 10).Value = 50
$ws.Cells.Item(533, 11).Value = 11000
$ws.Cells.Item(533, 12).Value = 11000
$ws.Cells.Item(533, 13).Value = 11000
$ws.Cells.Item(533, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(533, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(533, 16).Value = 5500
$ws.Cells.Item(533, 17).Value = 2
$ws.Cells.Item(533, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(533, 4).NumberFormat = $ws.Cells.Item(534, 4).NumberFormat
